$wb = $excel.ActiveWorkbook

# Both "展览" and "全部类型" sheets carry the same "想去人数" (F column)
# figures for these five events; bump them to the refreshed scrape values.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 650
    $ws.Range("F3").Value = 494
    $ws.Range("F4").Value = 35
    $ws.Range("F8").Value = 1969
    $ws.Range("F9").Value = 4070
}
